$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 116, shifting the existing rows 116:170 down to 117:171
# (this also pushes the sheet's used-range dimension from R170 to R171).
$ws.Rows.Item(116).Insert()

# Populate the freshly inserted row 116 with the new record.
$ws.Range("A116").Value = 10
$ws.Range("B116").Value = "Vega Modelo de Temuco"
$ws.Range("C116").Value = "La Araucanía"
$ws.Range("D116").Value2 = 44726
$ws.Range("E116").Value = 9
$ws.Range("F116").Value = 100114007
$ws.Range("G116").Value = "Jengibre"
$ws.Range("H116").Value = "Sin especificar"
$ws.Range("I116").Value = "Primera"
$ws.Range("J116").Value = 15
$ws.Range("K116").Value = 32500
$ws.Range("L116").Value = 32500
$ws.Range("M116").Value = 32500
$ws.Range("N116").Value = "$/caja 13 kilos"
$ws.Range("O116").Value = "Perú"
$ws.Range("P116").Value = 2500
$ws.Range("Q116").Value = 13
$ws.Range("R116").Value = "Hortaliza"
